$d = $word.ActiveDocument
$d.Content.Find.Execute("acitvly hindered out ability", $true, $false, $false, $false, $false,
                         $true, 1, $false, "actively hindered out ability", 2)
